# Gemini Facts.xlsx: append a spacecraft mass-budget block (rows 31-35)
# below the existing "Spacecraft" section, and restore the author's last
# view state (150% zoom, K12 selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A31").Value = "Launch Mass"
$ws.Range("B31").Value = "3850 kg"

$ws.Range("A32").Value = "Crew Size"
$ws.Range("B32").Value = 2

# Cells were originally authored in this interleaved order (not strictly
# row-by-row), which is what determines the order new strings land in the
# shared-string table -- reproduce it verbatim.
$ws.Range("B33").Value = "1982 kg"
$ws.Range("A34").Value = "Retrograde Module"
$ws.Range("A33").Value = "Re-entry Module"
$ws.Range("B34").Value = "591 kg"

$ws.Range("A35").Value = "Equiptment Module"
$ws.Range("B35").Value = "1278 kg"

# View-state changes recorded in the diff: zoom 110% -> 150%, selection
# moved from G23 to K12.
$excel.ActiveWindow.Zoom = 150
$ws.Range("K12").Select()
